$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 182854.55
$ws.Range("I28").Value = 182854.55
$ws.Range("K28").Value = 182854.55
$ws.Range("M28").Value = -182369.55
$ws.Range("H53").Value = 1151.9412
$ws.Range("I53").Value = 735.4
$ws.Range("J53").Value = 1747
$ws.Range("K53").Value = 735.4
$ws.Range("L53").Value = 1747
$ws.Range("M53").Value = -98.39999999999998
$ws.Range("N53").Value = -3021
$ws.Range("H106").Value = 2383403.8
$ws.Range("I106").Value = 2669116.2
$ws.Range("K106").Value = 2669116.2
$ws.Range("M106").Value = -2668485.2
$ws.Range("H129").Value = 1990
$ws.Range("I129").Value = 1430
$ws.Range("K129").Value = 4290
$ws.Range("M129").Value = 710
$ws.Range("H137").Value = 25402.084
$ws.Range("I137").Value = 32035.357
$ws.Range("J137").Value = 2185.625
$ws.Range("K137").Value = 96106.071
$ws.Range("L137").Value = 6556.875
$ws.Range("M137").Value = -93556.071
$ws.Range("N137").Value = -11656.875
$ws.Range("H138").Value = 2785.02
$ws.Range("I138").Value = 1297.5
$ws.Range("J138").Value = 3307.662
$ws.Range("K138").Value = 3892.5
$ws.Range("L138").Value = 9922.985999999999
$ws.Range("M138").Value = 1247.5
$ws.Range("N138").Value = -20202.986
$ws.Range("H141").Value = 14485.448
$ws.Range("I141").Value = 717.0714
$ws.Range("K141").Value = 2151.2142
$ws.Range("M141").Value = 3028.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25146642
$ws.Range("I32").Value = 27092174
$ws.Range("K32").Value = 27092174
$ws.Range("M32").Value = -27091887
$ws.Range("H61").Value = 4006.1667
$ws.Range("I61").Value = 3354.6843
$ws.Range("K61").Value = 3354.6843
$ws.Range("M61").Value = -3142.6843
$ws.Range("H74").Value = 2464.3
$ws.Range("I74").Value = 2337.6924
$ws.Range("K74").Value = 2337.6924
$ws.Range("M74").Value = -1463.6924
$ws.Range("H77").Value = 2464.3
$ws.Range("I77").Value = 2337.6924
$ws.Range("K77").Value = 11688.462
$ws.Range("M77").Value = -7320.462
$ws.Range("H132").Value = 252259.92
$ws.Range("I132").Value = 346252.72
$ws.Range("J132").Value = 4460.727
$ws.Range("K132").Value = 1038758.16
$ws.Range("L132").Value = 13382.181
$ws.Range("M132").Value = -1036228.16
$ws.Range("N132").Value = -18442.181
$ws.Range("H136").Value = 4006.1667
$ws.Range("I136").Value = 3354.6843
$ws.Range("K136").Value = 10064.0529
$ws.Range("M136").Value = -7514.052899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2673.75
$ws.Range("I99").Value = 2027.7142
$ws.Range("K99").Value = 2027.7142
$ws.Range("M99").Value = -529.7141999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4058.5
$ws.Range("I31").Value = 3119.5386
$ws.Range("J31").Value = 4872.2666
$ws.Range("K31").Value = 3119.5386
$ws.Range("L31").Value = 4872.2666
$ws.Range("M31").Value = -2824.5386
$ws.Range("N31").Value = -5462.2666
$ws.Range("H34").Value = 4058.5
$ws.Range("I34").Value = 3119.5386
$ws.Range("J34").Value = 4872.2666
$ws.Range("K34").Value = 3119.5386
$ws.Range("L34").Value = 4872.2666
$ws.Range("M34").Value = -2917.5386
$ws.Range("N34").Value = -5276.2666
$ws.Range("H58").Value = 3219.7097
$ws.Range("I58").Value = 2895.25
$ws.Range("K58").Value = 2895.25
$ws.Range("M58").Value = -2692.25
$ws.Range("H134").Value = 2536.8975
$ws.Range("I134").Value = 2445.3845
$ws.Range("J134").Value = 2719.923
$ws.Range("K134").Value = 7336.1535
$ws.Range("L134").Value = 8159.768999999999
$ws.Range("M134").Value = -4801.1535
$ws.Range("N134").Value = -13229.769
$ws.Range("H136").Value = 3219.7097
$ws.Range("I136").Value = 2895.25
$ws.Range("K136").Value = 8685.75
$ws.Range("M136").Value = -6135.75
$ws.Range("H141").Value = 538844.1
$ws.Range("J141").Value = 573048
$ws.Range("L141").Value = 573048
$ws.Range("N141").Value = -583408

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 85.333336
$ws.Range("I8").Value = 85.333336
$ws.Range("K8").Value = 256.000008
$ws.Range("M8").Value = -117.000008
$ws.Range("H11").Value = 750511.6
$ws.Range("I11").Value = 833682.3
$ws.Range("J11").Value = 500999.5
$ws.Range("K11").Value = 2501046.9
$ws.Range("L11").Value = 1502998.5
$ws.Range("M11").Value = -2500906.9
$ws.Range("N11").Value = -1503278.5
$ws.Range("H26").Value = 729.8333
$ws.Range("I26").Value = 729.8333
$ws.Range("K26").Value = 2189.4999
$ws.Range("M26").Value = -1901.4999
$ws.Range("H115").Value = 6000
$ws.Range("J115").Value = 7000
$ws.Range("L115").Value = 21000
$ws.Range("N115").Value = -23350
$ws.Range("H127").Value = 2236.375
$ws.Range("J127").Value = 2236.375
$ws.Range("L127").Value = 6709.125
$ws.Range("N127").Value = -16629.125
$ws.Range("H132").Value = 1138.4117
$ws.Range("I132").Value = 1089.7693
$ws.Range("K132").Value = 9807.923699999999
$ws.Range("M132").Value = -7277.923699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2690444.5
$ws.Range("I11").Value = 926333.3
$ws.Range("J11").Value = 3572500
$ws.Range("K11").Value = 926333.3
$ws.Range("L11").Value = 3572500
$ws.Range("M11").Value = -926194.3
$ws.Range("N11").Value = -3572778
$ws.Range("H122").Value = 1581
$ws.Range("I122").Value = 1318
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 3954
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1504
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12179.363
$ws.Range("I40").Value = 11248.375
$ws.Range("K40").Value = 11248.375
$ws.Range("M40").Value = -11112.375
$ws.Range("H100").Value = 2735.889
$ws.Range("I100").Value = 1783
$ws.Range("K100").Value = 1783
$ws.Range("M100").Value = -1242

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18782.197
$ws.Range("I132").Value = 20568.814
$ws.Range("K132").Value = 61706.442
$ws.Range("M132").Value = -59176.442
$ws.Range("H136").Value = 2527.3872
$ws.Range("I136").Value = 2137.5
$ws.Range("J136").Value = 3864.1428
$ws.Range("K136").Value = 6412.5
$ws.Range("L136").Value = 11592.4284
$ws.Range("M136").Value = -3862.5
$ws.Range("N136").Value = -16692.4284
